$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual cell values to reflect the refreshed stock-screener lists
$ws.Range("C2").Value = "NSE:ADVENZYMES"
$ws.Range("D2").Value = "NSE:DIXON"
$ws.Range("E2").Value = "NSE:ICICIPRULI"
$ws.Range("F2").Value = "NSE:HAL"
$ws.Range("B3").Value = "NSE:BDL"
$ws.Range("C3").Value = "NSE:ALBERTDAVD"
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "NSE:IDFCFIRSTB"
$ws.Range("B4").Value = "NSE:CRISIL"
$ws.Range("C4").Value = "NSE:BAJAJHIND"
$ws.Range("D4").ClearContents()
$ws.Range("B5").Value = "NSE:DCM"
$ws.Range("C5").Value = "NSE:BALKRISHNA"
$ws.Range("E5").ClearContents()
$ws.Range("B6").Value = "NSE:DCMNVL"
$ws.Range("C6").Value = "NSE:BTML"
$ws.Range("B7").Value = "NSE:EMUDHRA"
$ws.Range("C7").Value = "NSE:CCL"
$ws.Range("B8").Value = "NSE:HAL"
$ws.Range("C8").Value = "NSE:CINELINE"
$ws.Range("B9").Value = "NSE:HDFCNIFIT"
$ws.Range("C9").Value = "NSE:CONSOFINVT"
$ws.Range("B10").Value = "NSE:INFOMEDIA"
$ws.Range("C10").Value = "NSE:CSBBANK"
$ws.Range("B11").Value = "NSE:ITBEES"
$ws.Range("C11").Value = "NSE:DPWIRES"
$ws.Range("B12").Value = "NSE:MANAPPURAM"
$ws.Range("C12").Value = "NSE:FCSSOFT"
$ws.Range("B13").Value = "NSE:MVGJL"
$ws.Range("C13").Value = "NSE:GARFIBRES"
$ws.Range("B14").Value = "NSE:PSPPROJECT"
$ws.Range("C14").Value = "NSE:GENSOL"
$ws.Range("C15").Value = "NSE:HISARMETAL"
$ws.Range("C16").Value = "NSE:HPIL"
$ws.Range("C17").Value = "NSE:KBCGLOBAL"
$ws.Range("C18").Value = "NSE:KICL"
$ws.Range("C19").Value = "NSE:MAWANASUG"
$ws.Range("C20").Value = "NSE:MPSLTD"
$ws.Range("C21").Value = "NSE:MUTHOOTCAP"
$ws.Range("C22").Value = "NSE:NYKAA"
$ws.Range("C23").Value = "NSE:ONEPOINT"
$ws.Range("C24").Value = "NSE:PNC"
$ws.Range("C25").Value = "NSE:PRECWIRE"
$ws.Range("C26").Value = "NSE:PTCIL"
$ws.Range("C27").Value = "NSE:RAMCOSYS"
$ws.Range("C28").Value = "NSE:RBL"
$ws.Range("C29").Value = "NSE:SALASAR"

# Remove the trailing rows that no longer have data (rows 30-35)
$ws.Range("A30:F35").EntireRow.Delete()
